$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.178.86'
$ws.Range("E2").Value = '  +0.43%  '

$ws.Range("D3").Value = '1.637.82'
$ws.Range("E3").Value = '  -0.32%  '

$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("E5").Value = '  -0.46%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.518'
$ws.Range("E6").Value = '  +1.93%  '

$ws.Range("E7").Value = '  +0.26%  '

$ws.Range("E8").Value = '  -0.44%  '

$ws.Range("E9").Value = '  +0.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.13'
$ws.Range("E10").Value = '  -0.88%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0849'
$ws.Range("E11").Value = '  +0.30%  '

$ws.Range("D12").Value = '1.867.58'
$ws.Range("E12").Value = '  -0.22%  '

$ws.Range("D13").Value = '1.643.60'
$ws.Range("E13").Value = '  +0.15%  '

$ws.Range("E14").Value = '  +0.32%  '

$ws.Range("E15").Value = '  +0.94%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.44'
$ws.Range("E16").Value = '  -1.36%  '

$ws.Range("D17").Value = '27.185.08'
$ws.Range("E17").Value = '  +0.55%  '

$ws.Range("E18").Value = '  +0.41%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '217.17'
$ws.Range("E19").Value = '  -1.75%  '

$ws.Range("E20").Value = '  +0.14%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.86'
$ws.Range("E21").Value = '  +1.63%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.55'
$ws.Range("E22").Value = '  +4.51%  '

$ws.Range("E23").Value = '  -0.35%  '

$ws.Range("E24").Value = '  -1.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.56'
$ws.Range("E25").Value = '  +0.08%  '

$ws.Range("E26").Value = '  +0.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.38'
$ws.Range("E27").Value = '  -0.18%  '

$ws.Range("E28").Value = '  -0.33%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.68'
$ws.Range("E29").Value = '  -0.90%  '

$ws.Range("E30").Value = '  +0.20%  '

$ws.Range("E31").Value = '  -0.33%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.39'
$ws.Range("E32").Value = '  +1.21%  '

$ws.Range("E33").Value = '  +0.01%  '

$ws.Range("D34").Value = '1.301.94'
$ws.Range("E34").Value = '  +2.33%  '

$ws.Range("E35").Value = '  -0.36%  '

$ws.Range("E36").Value = '  +0.74%  '

$ws.Range("E37").Value = '  -0.64%  '

$ws.Range("E38").Value = '  +0.54%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.854'
$ws.Range("E39").Value = '  +1.19%  '

$ws.Range("E40").Value = '  +0.08%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.810'
$ws.Range("E41").Value = '  +0.04%  '

$ws.Range("E42").Value = '  +5.49%  '

$ws.Range("E43").Value = '  -0.55%  '

$ws.Range("D44").Value = '1.777.66'
$ws.Range("E44").Value = '  -0.22%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.44'
$ws.Range("E45").Value = '  -0.71%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.60'
$ws.Range("E46").Value = '  -2.26%  '

$ws.Range("E47").Value = '  -1.24%  '

$ws.Range("D48").Value = '0.0₆0106'
$ws.Range("E48").Value = '  +0.49%  '

$ws.Range("E49").Value = '  +0.41%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.62'
$ws.Range("E50").Value = '  -0.90%  '

$ws.Range("B51").Value = 'WEMIXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.758'
$ws.Range("E51").Value = '  +14.59%  '
